$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header "time_taken" — copy the header formatting (bold,
# bordered, centered) from the adjacent E1 header cell so F1 matches the
# other header cells' style, then set its text.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Per-row "time_taken" metadata values for rows 2-26 (unstyled, matching
# the rest of the data cells).
$timeTaken = @(
  "2021-10-05 10:50:56.747576",
  "2021-10-05 10:50:56.747588",
  "2021-10-05 10:50:56.747592",
  "2021-10-05 10:50:56.747594",
  "2021-10-05 10:50:56.747597",
  "2021-10-05 10:50:56.747600",
  "2021-10-05 10:50:56.747602",
  "2021-10-05 10:50:56.747605",
  "2021-10-05 10:50:56.747608",
  "2021-10-05 10:50:56.747610",
  "2021-10-05 10:50:56.747613",
  "2021-10-05 10:50:56.747616",
  "2021-10-05 10:50:56.747618",
  "2021-10-05 10:50:56.747621",
  "2021-10-05 10:50:56.747623",
  "2021-10-05 10:50:56.747626",
  "2021-10-05 10:50:56.747629",
  "2021-10-05 10:50:56.747631",
  "2021-10-05 10:50:56.747634",
  "2021-10-05 10:50:56.747636",
  "2021-10-05 10:50:56.747639",
  "2021-10-05 10:50:56.747641",
  "2021-10-05 10:50:56.747644",
  "2021-10-05 10:50:56.747647",
  "2021-10-05 10:50:56.747649"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
